{"js": "// Replace \"BRUNA PETRONI CEZARIO\" / \"Gerente de HSE Brasil\" (HSE manager\n// signature block) with \"LEONARDO SILVERIO FERREIRA\" / \"T\u00e9cnico(a) de\n// Seguran\u00e7a do Trabalho\", matching the formatting tweaks from the diff:\n//   - name paragraph mark (pPr/rPr) gains <w:b/><w:bCs/> and loses <w:lang/>\n//   - role paragraph run gains <w:lang w:val=\"pt-BR\"/>\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet nameParagraph = null;\nlet roleParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.trim();\n  if (text === \"BRUNA PETRONI CEZARIO\") {\n    nameParagraph = paragraphs.items[i];\n  } else if (text === \"Gerente de HSE Brasil\") {\n    roleParagraph = paragraphs.items[i];\n  }\n}\n\nif (!nameParagraph || !roleParagraph) {\n  throw new Error(\"Could not locate the BRUNA PETRONI CEZARIO / Gerente de HSE Brasil paragraphs\");\n}\n\n// Paragraph 1: the name line. Bold the paragraph mark itself (pPr/rPr)\n// and drop its <w:lang> entry, then swap the run text (run formatting\n// \u2014 bold/bCs/color \u2014 stays the same).\nconst nameRange = nameParagraph.getRange(\"Whole\");\nconst nameOoxml = `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n        <w:body>\n          <w:p w14:paraId=\"0F80B0E9\" w14:textId=\"77777777\" w:rsidR=\"001C081F\" w:rsidRPr=\"00123AA2\" w:rsidRDefault=\"001C081F\" w:rsidP=\"001C081F\">\n            <w:pPr>\n              <w:pStyle w:val=\"TableParagraph\"/>\n              <w:jc w:val=\"center\"/>\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n                <w:color w:val=\"000000\" w:themeColor=\"text1\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r w:rsidRPr=\"00123AA2\">\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n                <w:color w:val=\"000000\" w:themeColor=\"text1\"/>\n              </w:rPr>\n              <w:t>LEONARDO SILVERIO FERREIRA</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\nnameRange.insertOoxml(nameOoxml, \"Replace\");\nawait context.sync();\n\n// Paragraph 2: the role line. Add <w:lang w:val=\"pt-BR\"/> to the run and\n// swap the text; paragraph mark formatting is unchanged.\nconst roleRange = roleParagraph.getRange(\"Whole\");\nconst roleOoxml = `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n        <w:body>\n          <w:p w14:paraId=\"570F9CCE\" w14:textId=\"77777777\" w:rsidR=\"001C081F\" w:rsidRPr=\"00123AA2\" w:rsidRDefault=\"001C081F\" w:rsidP=\"001C081F\">\n            <w:pPr>\n              <w:pStyle w:val=\"TableParagraph\"/>\n              <w:jc w:val=\"center\"/>\n              <w:rPr>\n                <w:color w:val=\"000000\" w:themeColor=\"text1\"/>\n                <w:lang w:val=\"pt-BR\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r w:rsidRPr=\"00123AA2\">\n              <w:rPr>\n                <w:color w:val=\"000000\" w:themeColor=\"text1\"/>\n                <w:lang w:val=\"pt-BR\"/>\n              </w:rPr>\n              <w:t>T\u00e9cnico(a) de Seguran\u00e7a do Trabalho</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\nroleRange.insertOoxml(roleOoxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Replace \"BRUNA PETRONI CEZARIO\" / \"Gerente de HSE Brasil\" (HSE manager\n# signature block) with \"LEONARDO SILVERIO FERREIRA\" / \"T\u00e9cnico(a) de\n# Seguran\u00e7a do Trabalho\", matching the formatting tweaks from the diff:\n#   - name paragraph mark (pPr/rPr) gains <w:b/><w:bCs/> and loses <w:lang/>\n#   - role paragraph run gains <w:lang w:val=\"pt-BR\"/>\n\n$d = $word.ActiveDocument\n\n$nameParaIndex = $null\n$roleParaIndex = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text.Trim()\n    if ($t -eq \"BRUNA PETRONI CEZARIO\") {\n        $nameParaIndex = $i\n    } elseif ($t -eq \"Gerente de HSE Brasil\") {\n        $roleParaIndex = $i\n    }\n}\n\nif (-not $nameParaIndex -or -not $roleParaIndex) {\n    throw \"Could not locate the BRUNA PETRONI CEZARIO / Gerente de HSE Brasil paragraphs\"\n}\n\n# Paragraph 1: the name line. Bold the paragraph mark itself (pPr/rPr)\n# and drop its <w:lang> entry, then swap the run text (run formatting\n# -- bold/bCs/color -- stays the same).\n$nameRange = $d.Paragraphs($nameParaIndex).Range\n$nameOoxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n        <w:body>\n          <w:p w14:paraId=\"0F80B0E9\" w14:textId=\"77777777\" w:rsidR=\"001C081F\" w:rsidRPr=\"00123AA2\" w:rsidRDefault=\"001C081F\" w:rsidP=\"001C081F\">\n            <w:pPr>\n              <w:pStyle w:val=\"TableParagraph\"/>\n              <w:jc w:val=\"center\"/>\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n                <w:color w:val=\"000000\" w:themeColor=\"text1\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r w:rsidRPr=\"00123AA2\">\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n                <w:color w:val=\"000000\" w:themeColor=\"text1\"/>\n              </w:rPr>\n              <w:t>LEONARDO SILVERIO FERREIRA</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>'\n$nameRange.InsertXML($nameOoxml)\n\n# Paragraph 2: the role line. Add <w:lang w:val=\"pt-BR\"/> to the run and\n# swap the text; paragraph mark formatting is unchanged.\n$roleRange = $d.Paragraphs($roleParaIndex).Range\n$roleOoxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n        <w:body>\n          <w:p w14:paraId=\"570F9CCE\" w14:textId=\"77777777\" w:rsidR=\"001C081F\" w:rsidRPr=\"00123AA2\" w:rsidRDefault=\"001C081F\" w:rsidP=\"001C081F\">\n            <w:pPr>\n              <w:pStyle w:val=\"TableParagraph\"/>\n              <w:jc w:val=\"center\"/>\n              <w:rPr>\n                <w:color w:val=\"000000\" w:themeColor=\"text1\"/>\n                <w:lang w:val=\"pt-BR\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r w:rsidRPr=\"00123AA2\">\n              <w:rPr>\n                <w:color w:val=\"000000\" w:themeColor=\"text1\"/>\n                <w:lang w:val=\"pt-BR\"/>\n              </w:rPr>\n              <w:t>T\u00e9cnico(a) de Seguran\u00e7a do Trabalho</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>'\n$roleRange.InsertXML($roleOoxml)\n"}
